$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 74.75
$ws.Cells.Item(9, 9).Value = 74.75
$ws.Cells.Item(9, 11).Value = 74.75
$ws.Cells.Item(9, 13).Value = 94.25

$ws.Cells.Item(38, 8).Value = 374
$ws.Cells.Item(38, 9).Value = 61
$ws.Cells.Item(38, 11).Value = 183
$ws.Cells.Item(38, 13).Value = 189

$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = ""
$ws.Cells.Item(70, 14).Value = ""

$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = ""
$ws.Cells.Item(73, 14).Value = ""

$ws.Cells.Item(127, 8).Value = 3197
$ws.Cells.Item(127, 9).Value = 3197
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 9591
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 13).Value = -4631
$ws.Cells.Item(127, 14).Value = ""

$ws.Cells.Item(137, 8).Value = 8126.636
$ws.Cells.Item(137, 9).Value = 6200.2
$ws.Cells.Item(137, 11).Value = 18600.6
$ws.Cells.Item(137, 13).Value = -16050.6

$ws.Cells.Item(138, 8).Value = 9250

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 2214
$ws.Cells.Item(88, 9).Value = 500
$ws.Cells.Item(88, 10).Value = 2499.6667
$ws.Cells.Item(88, 11).Value = 500
$ws.Cells.Item(88, 12).Value = 2499.6667
$ws.Cells.Item(88, 13).Value = -94
$ws.Cells.Item(88, 14).Value = -3311.6667

$ws.Cells.Item(91, 8).Value = 2214
$ws.Cells.Item(91, 9).Value = 500
$ws.Cells.Item(91, 10).Value = 2499.6667
$ws.Cells.Item(91, 11).Value = 500
$ws.Cells.Item(91, 12).Value = 2499.6667
$ws.Cells.Item(91, 13).Value = 904
$ws.Cells.Item(91, 14).Value = -5307.6667

$ws.Cells.Item(122, 8).Value = 5000
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 14).Value = -19900

$ws.Cells.Item(124, 8).Value = 36927.93
$ws.Cells.Item(124, 10).Value = 36927.93
$ws.Cells.Item(124, 12).Value = 36927.93
$ws.Cells.Item(124, 14).Value = -46747.93

$ws.Cells.Item(132, 8).Value = 15398.1
$ws.Cells.Item(132, 9).Value = 13330.167
$ws.Cells.Item(132, 10).Value = 18500
$ws.Cells.Item(132, 11).Value = 39990.501
$ws.Cells.Item(132, 12).Value = 55500
$ws.Cells.Item(132, 13).Value = -37460.501
$ws.Cells.Item(132, 14).Value = -60560

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 30000
$ws.Cells.Item(26, 9).Value = 30000
$ws.Cells.Item(26, 11).Value = 30000
$ws.Cells.Item(26, 13).Value = -29708

$ws.Cells.Item(99, 8).Value = 3249.5
$ws.Cells.Item(99, 9).Value = 3249.5
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 3249.5
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -1751.5
$ws.Cells.Item(99, 14).Value = ""

$ws.Cells.Item(105, 8).Value = 2660.75
$ws.Cells.Item(105, 9).Value = 2660.75
$ws.Cells.Item(105, 11).Value = 2660.75
$ws.Cells.Item(105, 13).Value = -913.75

$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 14).Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6632.909
$ws.Cells.Item(31, 9).Value = 1428.6
$ws.Cells.Item(31, 10).Value = 10969.833
$ws.Cells.Item(31, 11).Value = 1428.6
$ws.Cells.Item(31, 12).Value = 10969.833
$ws.Cells.Item(31, 13).Value = -1133.6
$ws.Cells.Item(31, 14).Value = -11559.833

$ws.Cells.Item(34, 8).Value = 6632.909
$ws.Cells.Item(34, 9).Value = 1428.6
$ws.Cells.Item(34, 10).Value = 10969.833
$ws.Cells.Item(34, 11).Value = 1428.6
$ws.Cells.Item(34, 12).Value = 10969.833
$ws.Cells.Item(34, 13).Value = -1226.6
$ws.Cells.Item(34, 14).Value = -11373.833

$ws.Cells.Item(80, 8).Value = 36000
$ws.Cells.Item(80, 10).Value = 36000
$ws.Cells.Item(80, 12).Value = 36000
$ws.Cells.Item(80, 14).Value = -38246

$ws.Cells.Item(83, 8).Value = 36000
$ws.Cells.Item(83, 10).Value = 36000
$ws.Cells.Item(83, 12).Value = 108000
$ws.Cells.Item(83, 14).Value = -119232

$ws.Cells.Item(107, 8).Value = 4794.8
$ws.Cells.Item(107, 10).Value = 778
$ws.Cells.Item(107, 12).Value = 778
$ws.Cells.Item(107, 14).Value = -4618

$ws.Cells.Item(132, 8).Value = 10999.5
$ws.Cells.Item(132, 10).Value = 12000
$ws.Cells.Item(132, 12).Value = 36000
$ws.Cells.Item(132, 14).Value = -41060

$ws.Cells.Item(134, 8).Value = 10382.167
$ws.Cells.Item(134, 9).Value = 3899
$ws.Cells.Item(134, 10).Value = 13623.75
$ws.Cells.Item(134, 11).Value = 11697
$ws.Cells.Item(134, 12).Value = 40871.25
$ws.Cells.Item(134, 13).Value = -9162
$ws.Cells.Item(134, 14).Value = -45941.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 40.466667
$ws.Cells.Item(2, 9).Value = 48.5
$ws.Cells.Item(2, 10).Value = 8.333333
$ws.Cells.Item(2, 11).Value = 291
$ws.Cells.Item(2, 12).Value = 49.999998
$ws.Cells.Item(2, 13).Value = -178
$ws.Cells.Item(2, 14).Value = -275.999998

$ws.Cells.Item(12, 8).Value = 109.8
$ws.Cells.Item(12, 9).Value = 18.5
$ws.Cells.Item(12, 11).Value = 55.5
$ws.Cells.Item(12, 13).Value = 117.5

$ws.Cells.Item(38, 8).Value = 205
$ws.Cells.Item(38, 9).Value = 125
$ws.Cells.Item(38, 10).Value = 231.66667
$ws.Cells.Item(38, 11).Value = 375
$ws.Cells.Item(38, 12).Value = 695.00001
$ws.Cells.Item(38, 13).Value = -28
$ws.Cells.Item(38, 14).Value = -1389.00001

$ws.Cells.Item(119, 8).Value = 564.5
$ws.Cells.Item(119, 9).Value = 564.5
$ws.Cells.Item(119, 11).Value = 1693.5
$ws.Cells.Item(119, 13).Value = 3144.5

$ws.Cells.Item(131, 8).Value = 2341.25
$ws.Cells.Item(131, 10).Value = 2463.182
$ws.Cells.Item(131, 12).Value = 7389.545999999999
$ws.Cells.Item(131, 14).Value = -17469.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3391.1428
$ws.Cells.Item(102, 9).Value = 3391.1428
$ws.Cells.Item(102, 11).Value = 3391.1428
$ws.Cells.Item(102, 13).Value = -1769.1428

$ws.Cells.Item(122, 8).Value = 3625
$ws.Cells.Item(122, 9).Value = 3625
$ws.Cells.Item(122, 11).Value = 10875
$ws.Cells.Item(122, 13).Value = -8425

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 13).Value = ""

$ws.Cells.Item(61, 8).Value = 2000
$ws.Cells.Item(61, 9).Value = 2000
$ws.Cells.Item(61, 11).Value = 2000
$ws.Cells.Item(61, 13).Value = -1798

$ws.Cells.Item(113, 8).Value = 2000
$ws.Cells.Item(113, 9).Value = 2000
$ws.Cells.Item(113, 11).Value = 2000
$ws.Cells.Item(113, 13).Value = 170

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).Value = ""

$ws.Cells.Item(17, 8).Value = 4995
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 4995
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 4995
$ws.Cells.Item(17, 13).Value = ""
$ws.Cells.Item(17, 14).Value = -5339

$ws.Cells.Item(62, 8).Value = 2499
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 2499
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 2499
$ws.Cells.Item(62, 13).Value = ""
$ws.Cells.Item(62, 14).Value = -3747

$ws.Cells.Item(65, 8).Value = 2499
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 2499
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 12495
$ws.Cells.Item(65, 13).Value = ""
$ws.Cells.Item(65, 14).Value = -18735

$ws.Cells.Item(118, 8).Value = 65000
$ws.Cells.Item(118, 10).Value = 65000
$ws.Cells.Item(118, 12).Value = 65000
$ws.Cells.Item(118, 14).Value = -68314
